$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")
$ws.Activate()
